# "sync with main repo" - refresh the two Salesforce-sourced sample values on the
# "Routing Master" sheet (Item Number / Id) with the latest records, and re-fit the
# columns that hold them (their best-fit width grows slightly with the new text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

$ws.Range("B2").Value = "Pro-PEItem-NXQGE"
$ws.Range("D2").Value = "a345f000000uauIAAQ"

$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(4).ColumnWidth = 19.8
